$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Capture existing hyperlinks (row, column, target address, display text) ---
$links = @()
foreach ($hl in $ws.Hyperlinks) {
    $info = @{}
    $info.Row = $hl.Range.Row
    $info.Col = $hl.Range.Column
    $info.Address = $hl.Address
    $info.Display = $hl.TextToDisplay
    $links += $info
}

# --- 2. Delete row 2 (the snirelad61@gmail.com review) - shifts all rows below up by one ---
$ws.Rows(2).Delete()

# --- 3. Drop all (now stale) hyperlinks ---
$ws.Hyperlinks.Delete()

# --- 4. Re-create hyperlinks shifted to match the row deletion, skipping the one that lived on row 2 ---
foreach ($info in $links) {
    if ($info.Row -eq 2) {
        continue
    }
    $newRow = $info.Row - 1
    $addr = $ws.Cells.Item($newRow, $info.Col).Address()
    $ws.Hyperlinks.Add($ws.Range($addr), $info.Address, "", "", $info.Display)
}

# --- 5. Restore the selection to match the saved workbook state ---
$ws.Range("B5").Select()
